# Hotfix in start script: append two newly-received DTU log records
# (2022-10-24 13:43:45 and 2022-10-27 22:39:04) to every sheet of the
# Msg8705 workbook, mirroring the existing row layout on each tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Msg8705" (sheet1): columns A..H, existing data ends at row 213.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 214 - seed from row 213 (keeps the blank-string B column intact),
# then overwrite the cells that actually change.
$ws1.Range("A213:H213").Copy($ws1.Range("A214:H214"))
$ws1.Cells.Item(214, 1).Value = "2022-10-24 13:43:45"
$ws1.Cells.Item(214, 4).Value = 1666611832
$ws1.Cells.Item(214, 5).Value = 6
$ws1.Cells.Item(214, 6).Value = 24
$ws1.Cells.Item(214, 7).Value = 1
$ws1.Cells.Item(214, 8).Value = 255

# Row 215
$ws1.Range("A213:H213").Copy($ws1.Range("A215:H215"))
$ws1.Cells.Item(215, 1).Value = "2022-10-27 22:39:04"
$ws1.Cells.Item(215, 4).Value = 1666903154
$ws1.Cells.Item(215, 5).Value = 6
$ws1.Cells.Item(215, 6).Value = 24
$ws1.Cells.Item(215, 7).Value = 1
$ws1.Cells.Item(215, 8).Value = 255

# ---------------------------------------------------------------------
# Sheet "Msg8705_8" (sheet2): columns A..L, existing data ends at row 213.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Row 214
$ws2.Range("A213:L213").Copy($ws2.Range("A214:L214"))
$ws2.Cells.Item(214, 1).Value = "2022-10-24 13:43:45"
$ws2.Cells.Item(214, 3).Value = 527
$ws2.Cells.Item(214, 4).Value = 37122
$ws2.Cells.Item(214, 5).Value = 10
$ws2.Cells.Item(214, 6).Value = 160
$ws2.Cells.Item(214, 7).Value = 1
$ws2.Cells.Item(214, 11).Value = 21

# Row 215
$ws2.Range("A213:L213").Copy($ws2.Range("A215:L215"))
$ws2.Cells.Item(215, 1).Value = "2022-10-27 22:39:04"
$ws2.Cells.Item(215, 3).Value = 527
$ws2.Cells.Item(215, 4).Value = 37122
$ws2.Cells.Item(215, 5).Value = 10
$ws2.Cells.Item(215, 6).Value = 160
$ws2.Cells.Item(215, 7).Value = 1
$ws2.Cells.Item(215, 11).Value = 21

# ---------------------------------------------------------------------
# Sheet "Msg8705_11" (sheet3): columns A..H, existing data ends at row 637.
# None of the existing rows match the new numeric payloads, so seed the
# blank-string B column from row 637 and overwrite every value column.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$row = 638
$stamp231 = "2022-10-24 13:43:45"
$stamp232 = "2022-10-27 22:39:04"

$payloads = @(
    @(1914729089, 10012, 269627400, 256, 10752, 8192),
    @(1914726672, 10012, 269627400, 256, 10752, 8192),
    @(1914728742, 10012, 269627400, 256, 10752, 8192),
    @(-2145376735, 10016, 269627393, 256, 10752, 8192),
    @(-2145376224, 10016, 269627393, 256, 10752, 8192),
    @(-2145374071, 10016, 269627393, 256, 10752, 8192)
)

foreach ($stamp in @($stamp231, $stamp232)) {
    foreach ($payload in $payloads) {
        $ws3.Range("A637:H637").Copy($ws3.Range("A" + $row + ":H" + $row))
        $ws3.Cells.Item($row, 1).Value = $stamp
        $ws3.Cells.Item($row, 3).Value = $payload[0]
        $ws3.Cells.Item($row, 4).Value = $payload[1]
        $ws3.Cells.Item($row, 5).Value = $payload[2]
        $ws3.Cells.Item($row, 6).Value = $payload[3]
        $ws3.Cells.Item($row, 7).Value = $payload[4]
        $ws3.Cells.Item($row, 8).Value = $payload[5]
        $row = $row + 1
    }
}
